$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion rates inside the A1 text block ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$rngA1 = $wsHoja1.Range("A1")
$text = $rngA1.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 7.58 = 30579.92 pesos"), "1000 Bs = 7.59 = 30678.11 pesos"
$text = $text -replace [regex]::Escape("30579.92 pesos = 7.54 = 941.97 Bs"), "30678.11 pesos = 7.56 = 960.35 Bs"
$rngA1.Value2 = $text

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 131.69
$wsTasas.Range("O10").Value = 4040
$wsTasas.Range("N12").Value = 4057
$wsTasas.Range("O12").Value = 127
